# Update "想去人数" (want-to-go count, column F) figures to the freshly
# scraped numbers, and flip the COMICUP 2024SP ticket tier (column G) to
# "已售罄" (sold out) now that it no longer has a numeric minimum price.
# This mirrors the gh-pages data refresh generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) -----------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$sheet1Updates = @(
    @{ Row = 2;  F = 8392 },
    @{ Row = 3;  F = 36531 },   # COMICUP 2024SP -> also sells out (see below)
    @{ Row = 6;  F = 747 },
    @{ Row = 8;  F = 154 },
    @{ Row = 9;  F = 454 },
    @{ Row = 10; F = 822 },
    @{ Row = 11; F = 75 },
    @{ Row = 12; F = 653 },
    @{ Row = 13; F = 487 },
    @{ Row = 14; F = 29 },
    @{ Row = 15; F = 601 },
    @{ Row = 16; F = 169 },
    @{ Row = 17; F = 446 },
    @{ Row = 18; F = 432 },
    @{ Row = 19; F = 1137 },
    @{ Row = 21; F = 772 },
    @{ Row = 22; F = 2432 },
    @{ Row = 23; F = 925 },
    @{ Row = 24; F = 529 },
    @{ Row = 25; F = 88 },
    @{ Row = 26; F = 1124 },
    @{ Row = 28; F = 708 },
    @{ Row = 29; F = 29 },
    @{ Row = 30; F = 1120 }
)

foreach ($u in $sheet1Updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.F
}
$ws1.Cells.Item(3, 7).Value = "已售罄"

# --- 演出 (Performances) -----------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(5, 6).Value = 324

# --- 本地生活 (Local life) ----------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 589

# --- 全部类型 (All types, aggregates the other three sheets) -----------
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @(
    @{ Row = 2;  F = 589 },
    @{ Row = 3;  F = 8392 },
    @{ Row = 5;  F = 36531 },   # COMICUP 2024SP -> also sells out (see below)
    @{ Row = 8;  F = 747 },
    @{ Row = 11; F = 154 },
    @{ Row = 12; F = 454 },
    @{ Row = 14; F = 324 },
    @{ Row = 16; F = 822 },
    @{ Row = 17; F = 75 },
    @{ Row = 18; F = 653 },
    @{ Row = 19; F = 487 },
    @{ Row = 21; F = 29 },
    @{ Row = 26; F = 601 },
    @{ Row = 27; F = 169 },
    @{ Row = 28; F = 446 },
    @{ Row = 29; F = 432 },
    @{ Row = 30; F = 1137 },
    @{ Row = 32; F = 772 },
    @{ Row = 33; F = 2432 },
    @{ Row = 34; F = 925 },
    @{ Row = 35; F = 529 },
    @{ Row = 36; F = 88 },
    @{ Row = 37; F = 1124 },
    @{ Row = 40; F = 708 },
    @{ Row = 41; F = 29 },
    @{ Row = 42; F = 1120 }
)

foreach ($u in $sheet4Updates) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.F
}
$ws4.Cells.Item(5, 7).Value = "已售罄"
